$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.218632340431213
$ws.Range("B1").Value = 6.444764614105225
$ws.Range("C1").Value = 3.698637247085571
$ws.Range("D1").Value = 1.655169606208801
$ws.Range("E1").Value = 1.16774594783783
